$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.049.03"
$ws.Range("E2").Value = "  -3.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.642.27"
$ws.Range("E3").Value = "  -5.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.27"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.77"
$ws.Range("E6").Value = "  +5.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.640.63"
$ws.Range("E7").Value = "  -5.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.626"
$ws.Range("E8").Value = "  -5.77%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.709"
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("E11").Value = "  -8.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.02"
$ws.Range("E12").Value = "  +3.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000289"
$ws.Range("E13").Value = "  -10.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.60"
$ws.Range("E14").Value = "  -5.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.232.00"
$ws.Range("E15").Value = "  -5.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.651.18"
$ws.Range("E16").Value = "  -5.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.21"
$ws.Range("E17").Value = "  -9.22%  "
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.11"
$ws.Range("E19").Value = "  -6.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.64"
$ws.Range("E20").Value = "  -8.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.921.58"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "407.22"
$ws.Range("E22").Value = "  -6.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.54"
$ws.Range("E23").Value = "  -4.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.01"
$ws.Range("E24").Value = "  -6.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.98"
$ws.Range("E25").Value = "  -8.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.59"
$ws.Range("E26").Value = "  -8.71%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.88"
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.70"
$ws.Range("E28").Value = "  -6.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.05"
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.43"
$ws.Range("E30").Value = "  -9.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.41"
$ws.Range("E31").Value = "  -7.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.14"
$ws.Range("E32").Value = "  -13.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.26"
$ws.Range("E33").Value = "  -8.96%  "
$ws.Range("E34").Value = "  -7.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "64.42"
$ws.Range("E35").Value = "  -6.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "594.68"
$ws.Range("E36").Value = "  -6.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.44"
$ws.Range("E37").Value = "  -11.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0874"
$ws.Range("E38").Value = "  -10.23%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.394"
$ws.Range("E40").Value = "  -8.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  -6.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.00"
$ws.Range("E43").Value = "  -6.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.68"
$ws.Range("E44").Value = "  -7.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0435"
$ws.Range("E45").Value = "  -7.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.80"
$ws.Range("E46").Value = "  -12.10%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.71"
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.134"
$ws.Range("E48").Value = "  -6.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.92"
$ws.Range("E49").Value = "  -10.31%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.700.14"
$ws.Range("E50").Value = "  -6.64%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.09"
$ws.Range("E51").Value = "  -6.55%  "
